# remove comments from listnodes
$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "Sheet3"
$wsOld = $wb.Worksheets.Item("Sheet1")
$wsOld.Name = "Sheet3"

# Sheet2: update "List 1"/"Liste 1" entries to "List 2"/"Liste 2"
$wsSheet2 = $wb.Worksheets.Item("Sheet2")
$wsSheet2.Range("B2").Value = "List 2"
$wsSheet2.Range("C2").Value = "Liste 2"
$wsSheet2.Range("B3").Value = "List 2"
$wsSheet2.Range("C3").Value = "Liste 2"

# Update selection on Sheet2
$wsSheet2.Range("C21").Select()

# Sheet3 (formerly Sheet1): update "List 2" entries to "List 3"
$wsSheet3 = $wb.Worksheets.Item("Sheet3")
$wsSheet3.Range("A2").Value = "List 3"
$wsSheet3.Range("A3").Value = "List 3"
$wsSheet3.Range("A4").Value = "List 3"
$wsSheet3.Range("A5").Value = "List 3"
$wsSheet3.Range("A6").Value = "List 3"

# Update selection on Sheet3
$wsSheet3.Range("B16").Select()

# Re-activate Sheet2 as the active tab
$wsSheet2.Activate()
